# "Generate Report for Archive"
#
# Refresh the localization status report:
#   - flip the in-flight status label from "Ready for handoff" to
#     "In Translation" everywhere it appears (Overview!E2/F2,
#     zh-cn!C2, de-de!C2 all share that string)
#   - the new, shorter status text no longer needs as much horizontal
#     room, so the status columns are narrowed to fit it (Overview!E:F,
#     zh-cn!C, de-de!C)

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Narrowed status-column width (characters), matching the new text.
$statusColWidth = 12.42

# --- Overview sheet: status columns are E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E1:F1").ColumnWidth = $statusColWidth

# --- zh-cn sheet: status column is C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = $statusColWidth

# --- de-de sheet: status column is C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = $statusColWidth
